$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete rows 34-35 (SC 193 duplicate-shift source rows no longer needed;
# row "RM 232" is removed from the data set, shifting subsequent rows up,
# so the sheet shrinks from 35 to 33 data+header rows).
$ws.Range("A34:F35").Delete()

# Apply the corrected / re-imputed values for the remaining rows/cells.
$ws.Range("C2").Value = 14.9
$ws.Range("C3").Value = ""
$ws.Range("F3").Value = 17.64
$ws.Range("C4").Value = ""
$ws.Range("D5").Value = ""
$ws.Range("E6").Value = -5.7
$ws.Range("D8").Value = -13.9
$ws.Range("D10").Value = -14.7
$ws.Range("F10").Value = ""
$ws.Range("C11").Value = 11.4
$ws.Range("E11").Value = -7.9
$ws.Range("F11").Value = ""
$ws.Range("D12").Value = ""
$ws.Range("E12").Value = ""
$ws.Range("F12").Value = 17.45
$ws.Range("C13").Value = ""
$ws.Range("E13").Value = -5.3
$ws.Range("D15").Value = -15.2
$ws.Range("F16").Value = 17.34
$ws.Range("E17").Value = ""
$ws.Range("F17").Value = 17.78
$ws.Range("D18").Value = ""
$ws.Range("E18").Value = -8.5
$ws.Range("D19").Value = ""
$ws.Range("E19").Value = ""
$ws.Range("C21").Value = 12.7
$ws.Range("E24").Value = ""
$ws.Range("F24").Value = ""
$ws.Range("C25").Value = ""
$ws.Range("D25").Value = -15.5
$ws.Range("E25").Value = -7.1
$ws.Range("F25").Value = ""
$ws.Range("A26").Value = "SC 5"
$ws.Range("C26").Value = 10.8
$ws.Range("D26").Value = -13.8
$ws.Range("E26").Value = -5
$ws.Range("F26").Value = 17.38
$ws.Range("A27").Value = "SC 101"
$ws.Range("C27").Value = 10
$ws.Range("D27").Value = -14.6
$ws.Range("E27").Value = -10
$ws.Range("F27").Value = 17
$ws.Range("A28").Value = "SC 105"
$ws.Range("B28").Value = -19.6
$ws.Range("C28").Value = 11.1
$ws.Range("D28").Value = -13.7
$ws.Range("E28").Value = -5.9
$ws.Range("F28").Value = ""
$ws.Range("A29").Value = "SC 119"
$ws.Range("C29").Value = 11.2
$ws.Range("E29").Value = -6.8
$ws.Range("F29").Value = 18.06
$ws.Range("A30").Value = "SC 120"
$ws.Range("B30").Value = -19.7
$ws.Range("C30").Value = 11.4
$ws.Range("D30").Value = -13.6
$ws.Range("E30").Value = -5.7
$ws.Range("F30").Value = 16.89
$ws.Range("A31").Value = "SC 132"
$ws.Range("B31").Value = -18.8
$ws.Range("C31").Value = 15.3
$ws.Range("D31").Value = -13.7
$ws.Range("E31").Value = ""
$ws.Range("F31").Value = 17.18
$ws.Range("A32").Value = "SC 193"
$ws.Range("B32").Value = -19.9
$ws.Range("C32").Value = 10.5
$ws.Range("D32").Value = -14.7
$ws.Range("E32").Value = ""
$ws.Range("F32").Value = 17.39
$ws.Range("A33").Value = "SC 232"
$ws.Range("B33").Value = -19.5
$ws.Range("C33").Value = 10.4
$ws.Range("D33").Value = ""
$ws.Range("E33").Value = -10.7
$ws.Range("F33").Value = 17.53
